# Fruta / hortaliza, semanal
# Weekly refresh of the "Femacal de La Calera - Uva" dataset: three new
# price rows are inserted at the top of the data block (row 1116), which
# pushes the existing rows down by three positions (old 1116 -> new 1119,
# ..., old 1169 -> new 1172). The sheet dimension grows from A1:T1169 to
# A1:T1172.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 1116; everything that used to
# live at rows 1116-1169 shifts down to rows 1119-1172.
$ws.Rows("1116:1118").Insert()

# New row 1116: Red Globe
$ws.Cells.Item(1116, 1).Value = 3
$ws.Cells.Item(1116, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1116, 3).Value = "Coquimbo"
$ws.Cells.Item(1116, 4).Value = 45041
$ws.Cells.Item(1116, 5).Value = 5
$ws.Cells.Item(1116, 6).Value = "Fruta"
$ws.Cells.Item(1116, 7).Value = 100109
$ws.Cells.Item(1116, 8).Value = "Uva"
$ws.Cells.Item(1116, 9).Value = 100109001
$ws.Cells.Item(1116, 10).Value = "Uva"
$ws.Cells.Item(1116, 11).Value = "Red Globe"
$ws.Cells.Item(1116, 12).Value = "Primera"
$ws.Cells.Item(1116, 13).Value = 56
$ws.Cells.Item(1116, 14).Value = 7000
$ws.Cells.Item(1116, 15).Value = 7000
$ws.Cells.Item(1116, 16).Value = 7000
$ws.Cells.Item(1116, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(1116, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(1116, 19).Value = 700
$ws.Cells.Item(1116, 20).Value = 10

# New row 1117: Red Globe
$ws.Cells.Item(1117, 1).Value = 3
$ws.Cells.Item(1117, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1117, 3).Value = "Coquimbo"
$ws.Cells.Item(1117, 4).Value = 45041
$ws.Cells.Item(1117, 5).Value = 5
$ws.Cells.Item(1117, 6).Value = "Fruta"
$ws.Cells.Item(1117, 7).Value = 100109
$ws.Cells.Item(1117, 8).Value = "Uva"
$ws.Cells.Item(1117, 9).Value = 100109001
$ws.Cells.Item(1117, 10).Value = "Uva"
$ws.Cells.Item(1117, 11).Value = "Red Globe"
$ws.Cells.Item(1117, 12).Value = "Primera"
$ws.Cells.Item(1117, 13).Value = 70
$ws.Cells.Item(1117, 14).Value = 10000
$ws.Cells.Item(1117, 15).Value = 10000
$ws.Cells.Item(1117, 16).Value = 10000
$ws.Cells.Item(1117, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(1117, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(1117, 19).Value = 556
$ws.Cells.Item(1117, 20).Value = 18

# New row 1118: Thompson seedless
$ws.Cells.Item(1118, 1).Value = 3
$ws.Cells.Item(1118, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1118, 3).Value = "Coquimbo"
$ws.Cells.Item(1118, 4).Value = 45041
$ws.Cells.Item(1118, 5).Value = 5
$ws.Cells.Item(1118, 6).Value = "Fruta"
$ws.Cells.Item(1118, 7).Value = 100109
$ws.Cells.Item(1118, 8).Value = "Uva"
$ws.Cells.Item(1118, 9).Value = 100109001
$ws.Cells.Item(1118, 10).Value = "Uva"
$ws.Cells.Item(1118, 11).Value = "Thompson seedless"
$ws.Cells.Item(1118, 12).Value = "Primera"
$ws.Cells.Item(1118, 13).Value = 67
$ws.Cells.Item(1118, 14).Value = 8000
$ws.Cells.Item(1118, 15).Value = 8000
$ws.Cells.Item(1118, 16).Value = 8000
$ws.Cells.Item(1118, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(1118, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(1118, 19).Value = 800
$ws.Cells.Item(1118, 20).Value = 10
